# Update column G ("K") values for rows 2-22 on the active worksheet.
# These correspond to the regenerated K values (replacing the old Strike# values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 2
    4  = 3
    5  = 2
    6  = 2
    7  = 0
    8  = 1
    9  = 7
    10 = 3
    11 = 2
    12 = 7
    13 = 3
    14 = 2
    15 = 1
    16 = 2
    17 = 2
    18 = 1
    19 = 2
    20 = 2
    21 = 1
    22 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
